$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)
$ws2 = $wb.Worksheets.Item(2)

# --- Sheet1 (Tabelle1): content fixes ---
# Add new "Specification:" / "Estimated data" row labels
$ws1.Range("A34").Value = "Specification:"
$ws1.Range("B34").Value = "Estimated data"

# Match the new label's formatting to the existing similarly-styled label below it
$ws1.Range("A35").Copy()
$ws1.Range("A34").PasteSpecial(-4122)
$ws1.Range("B35").Copy()
$ws1.Range("B34").PasteSpecial(-4122)

# "Datenquelle:" -> "Source:"
$ws1.Range("A35").Value = "Source:"

# Fix typo'd data range for 2015 row (224 - 226 -> 224 - 266, matches Tabelle2 computed value)
$ws1.Range("B10").Value = "224 - 266"

# --- Chart: axis title translation ---
$co = $ws1.ChartObjects(1)
$chart = $co.Chart
$catAxis = $chart.Axes(1)
$catAxis.AxisTitle.Text = "Year"

# --- sheetView selection/navigation updates ---
$ws2.Range("D8").Select()
$ws1.Activate()
$ws1.Range("K30").Select()
